# Loan RBI, Variable Instalments
#
# 1. A new (blank-header) column is inserted on the "Repayment schedule"
#    sheet between "Paid" (K/L) and "In Advance" (old M, now N). The
#    existing "Due" figure moves from 10100 down to 100 and a new
#    "In Advance" style figure of 10000 is recorded.
# 2. The active/selected sheet+cell moves from "Prepay Loan" to the
#    "Repayment schedule" sheet, with the selection parked on I10.

$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the old column M (13th column), pushing
# "In Advance"/"Late"/heading/"Outstanding" one column to the right.
$ws.Columns.Item(13).Insert()

# Widen the "Paid" column (L, now still column 12) and the newly inserted
# column (M, column 13) to match the "Due" column's width.
$ws.Columns.Item(12).ColumnWidth = 5.71
$ws.Columns.Item(13).ColumnWidth = 5.71

# Update the repayment-schedule figures for row 3 (the single instalment row).
$ws.Cells.Item(3, 11).Value = 100    # K3: Due 10100 -> 100
$ws.Cells.Item(3, 14).Value = 10000  # N3 (old M3 position shifted): 0 -> 10000

# This sheet becomes the active tab/sheet with I10 selected.
$ws.Activate() | Out-Null
$ws.Range("I10").Select() | Out-Null
